# Update cryptos list data (prices and volume %) per commit
# Applies text-preserving updates to columns B-E for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '28.463.15'
$cell = $ws.Range('E2')
$cell.NumberFormat = '@'
$cell.Value = '  +0.10%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '1.819.76'
$cell = $ws.Range('E3')
$cell.NumberFormat = '@'
$cell.Value = '  -0.42%  '
$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '1.001'
$cell = $ws.Range('E4')
$cell.NumberFormat = '@'
$cell.Value = '  +0.00%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '315.26'
$cell = $ws.Range('E5')
$cell.NumberFormat = '@'
$cell.Value = '  -0.48%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '1.001'
$cell = $ws.Range('E6')
$cell.NumberFormat = '@'
$cell.Value = '  +0.02%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.5103'
$cell = $ws.Range('E7')
$cell.NumberFormat = '@'
$cell.Value = '  -4.64%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.3965'
$cell = $ws.Range('E8')
$cell.NumberFormat = '@'
$cell.Value = '  -1.53%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.08209'
$cell = $ws.Range('E9')
$cell.NumberFormat = '@'
$cell.Value = '  +6.75%  '
$cell = $ws.Range('B10')
$cell.NumberFormat = '@'
$cell.Value = 'Polygon'
$cell = $ws.Range('C10')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '1.110'
$cell = $ws.Range('E10')
$cell.NumberFormat = '@'
$cell.Value = '  -0.36%  '
$cell = $ws.Range('B11')
$cell.NumberFormat = '@'
$cell.Value = 'OKB'
$cell = $ws.Range('C11')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '41.57'
$cell = $ws.Range('E11')
$cell.NumberFormat = '@'
$cell.Value = '  -0.66%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '21.04'
$cell = $ws.Range('E12')
$cell.NumberFormat = '@'
$cell.Value = '  +0.27%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '6.290'
$cell = $ws.Range('E13')
$cell.NumberFormat = '@'
$cell.Value = '  -0.55%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell = $ws.Range('E14')
$cell.NumberFormat = '@'
$cell.Value = '  -0.01%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '7.523'
$cell = $ws.Range('E15')
$cell.NumberFormat = '@'
$cell.Value = '  -1.37%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '1.813.81'
$cell = $ws.Range('E16')
$cell.NumberFormat = '@'
$cell.Value = '  -0.86%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '0.00001138'
$cell = $ws.Range('E17')
$cell.NumberFormat = '@'
$cell.Value = '  +5.65%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '92.59'
$cell = $ws.Range('E18')
$cell.NumberFormat = '@'
$cell.Value = '  +3.18%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '0.06641'
$cell = $ws.Range('E19')
$cell.NumberFormat = '@'
$cell.Value = '  +0.80%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '17.72'
$cell = $ws.Range('E20')
$cell.NumberFormat = '@'
$cell.Value = '  -0.04%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell = $ws.Range('E21')
$cell.NumberFormat = '@'
$cell.Value = '  +0.05%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '6.108'
$cell = $ws.Range('E22')
$cell.NumberFormat = '@'
$cell.Value = '  +0.39%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '28.475.34'
$cell = $ws.Range('E23')
$cell.NumberFormat = '@'
$cell.Value = '  +0.12%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '11.38'
$cell = $ws.Range('E24')
$cell.NumberFormat = '@'
$cell.Value = '  +1.83%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '2.271'
$cell = $ws.Range('E25')
$cell.NumberFormat = '@'
$cell.Value = '  +1.91%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '21.28'
$cell = $ws.Range('E26')
$cell.NumberFormat = '@'
$cell.Value = '  +2.62%  '
$cell = $ws.Range('B27')
$cell.NumberFormat = '@'
$cell.Value = 'Monero'
$cell = $ws.Range('C27')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '155.46'
$cell = $ws.Range('E27')
$cell.NumberFormat = '@'
$cell.Value = '  -1.45%  '
$cell = $ws.Range('B28')
$cell.NumberFormat = '@'
$cell.Value = 'WrappedliquidstakedEther2.0'
$cell = $ws.Range('C28')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '2.024.75'
$cell = $ws.Range('E28')
$cell.NumberFormat = '@'
$cell.Value = '  -0.78%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '2.412'
$cell = $ws.Range('E29')
$cell.NumberFormat = '@'
$cell.Value = '  -2.17%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '125.68'
$cell = $ws.Range('E30')
$cell.NumberFormat = '@'
$cell.Value = '  +1.20%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '1.111'
$cell = $ws.Range('E31')
$cell.NumberFormat = '@'
$cell.Value = '  -1.25%  '
$cell = $ws.Range('E32')
$cell.NumberFormat = '@'
$cell.Value = '  -1.26%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '5.783'
$cell = $ws.Range('E33')
$cell.NumberFormat = '@'
$cell.Value = '  +1.85%  '
$cell = $ws.Range('E34')
$cell.NumberFormat = '@'
$cell.Value = '  +0.24%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '0.07063'
$cell = $ws.Range('E35')
$cell.NumberFormat = '@'
$cell.Value = '  -5.82%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.2224'
$cell = $ws.Range('E36')
$cell.NumberFormat = '@'
$cell.Value = '  -0.79%  '
$cell = $ws.Range('E37')
$cell.NumberFormat = '@'
$cell.Value = '  -0.32%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '5.228'
$cell = $ws.Range('E38')
$cell.NumberFormat = '@'
$cell.Value = '  -0.06%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '8.856'
$cell = $ws.Range('E39')
$cell.NumberFormat = '@'
$cell.Value = '  -0.39%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.6295'
$cell = $ws.Range('E40')
$cell.NumberFormat = '@'
$cell.Value = '  +0.33%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '11.28'
$cell = $ws.Range('E41')
$cell.NumberFormat = '@'
$cell.Value = '  -0.40%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '1.176'
$cell = $ws.Range('E42')
$cell.NumberFormat = '@'
$cell.Value = '  -0.13%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '1.000'
$cell = $ws.Range('E43')
$cell.NumberFormat = '@'
$cell.Value = '  +0.11%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '1.405'
$cell = $ws.Range('E44')
$cell.NumberFormat = '@'
$cell.Value = '  +0.81%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '13.54'
$cell = $ws.Range('E45')
$cell.NumberFormat = '@'
$cell.Value = '  -0.20%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '3.743'
$cell = $ws.Range('E46')
$cell.NumberFormat = '@'
$cell.Value = '  +1.18%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '0.5919'
$cell = $ws.Range('E47')
$cell.NumberFormat = '@'
$cell.Value = '  +1.08%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '125.33'
$cell = $ws.Range('E48')
$cell.NumberFormat = '@'
$cell.Value = '  +0.32%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '1.981'
$cell = $ws.Range('E49')
$cell.NumberFormat = '@'
$cell.Value = '  -1.25%  '
$cell = $ws.Range('E50')
$cell.NumberFormat = '@'
$cell.Value = '  -1.70%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '0.06893'
$cell = $ws.Range('E51')
$cell.NumberFormat = '@'
$cell.Value = '  -0.02%  '
